# Update countries & provincias Spain
# Refreshes the COVID country-stats table ("Pais" sheet) with newer figures
# and bumps the "last updated" timestamp. Also swaps the Suiza/Uzbekistan
# rows since Suiza's updated total now outranks Uzbekistan's.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" banner (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 13:57"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 6588825
$ws.Range("C4").Value = 662
$ws.Range("D4").Value = 3880706
$ws.Range("E4").Value = 2511774
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 196345

# --- India (row 5) ---
$ws.Range("B5").Value = 4568770
$ws.Range("C5").Value = 9045
$ws.Range("D5").Value = 3544794
$ws.Range("E5").Value = 947628
$ws.Range("G5").Value = 44
$ws.Range("H5").Value = 76348

# --- Catar (row 31) ---
$ws.Range("B31").Value = 121287
$ws.Range("C31").Value = 235
$ws.Range("D31").Value = 118199
$ws.Range("E31").Value = 2883

# --- Nepal (row 58) ---
$ws.Range("B58").Value = 51919
$ws.Range("C58").Value = 1454
$ws.Range("D58").Value = 36672
$ws.Range("E58").Value = 14925
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 322

# --- Suiza / Uzbekistan swap rows (60/61) ---
# Suiza's refreshed total (46239) now exceeds Uzbekistan's (45927), so Suiza
# moves up to row 60 and Uzbekistan drops to row 61.
$ws.Range("A60").Value = "Suiza"
$ws.Range("B60").Value = 46239
$ws.Range("C60").Value = 528
$ws.Range("D60").Value = 38100
$ws.Range("E60").Value = 6119
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 2020

$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 45927
$ws.Range("C61").Value = 454
$ws.Range("D61").Value = 42555
$ws.Range("E61").Value = 2997
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 375

# --- Estado de Palestina (row 73) ---
$ws.Range("B73").Value = 29256
$ws.Range("C73").Value = 592
$ws.Range("D73").Value = 19788
$ws.Range("E73").Value = 9264
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 204

# --- Madagascar (row 85) ---
$ws.Range("B85").Value = 15669
$ws.Range("C85").Value = 45
$ws.Range("D85").Value = 14320
$ws.Range("E85").Value = 1140
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 209

# --- Senegal (row 87) ---
$ws.Range("B87").Value = 14193
$ws.Range("C87").Value = 43
$ws.Range("D87").Value = 10350
$ws.Range("E87").Value = 3550

# --- Consejo Danes para los Refugiados (row 95) ---
$ws.Range("B95").Value = 10361
$ws.Range("C95").Value = 18
$ws.Range("D95").Value = 9622
$ws.Range("E95").Value = 477

# --- Gambia (row 131) ---
$ws.Range("B131").Value = 3362
$ws.Range("C131").Value = 32
$ws.Range("D131").Value = 1582
$ws.Range("E131").Value = 1680

# --- Vietnam (row 165) ---
$ws.Range("B165").Value = 1060
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 902
$ws.Range("E165").Value = 123
